# Commit: "Contact info changed in test case 520"
#
# The feedback-channel contact info on the CodeSchemes sheet (row 2,
# columns X/Y/Z = FEEDBACK_CHANNEL_FI/EN/SW-UG) changes from the
# "Ylapitajan/Yllapitajan yhteystiedot" (administrator contact) wording to
# "Aineiston palauteosoite" (data feedback address) wording, and the
# sheet's view/selection moves from X2 to Z2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CodeSchemes")

$ws.Range("X2").Value = "Aineiston palauteosoite_fi"
$ws.Range("Y2").Value = "Aineiston palauteosoite_en"
$ws.Range("Z2").Value = "Aineiston palauteosoite_sw-UG"

# Reflect the new active cell / scroll position recorded in the sheet view.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 21
$ws.Range("Z2").Select()
